$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")
$ws.Range("A1").Value = "Categories"
$ws.Range("A2").Value = "Laptops"
$ws.Range("A3").Value = "Accessories"

$rng = $ws.Range("A2:A3")
$rng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$rng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$rng.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$rng.Locked = $false
